$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.745.38'
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("D3").Value = '3.579.15'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''590.26'
$ws.Range("E5").Value = '  +2.73%  '
$ws.Range("D6").Value = '''187.88'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '3.568.85'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''0.622'
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '''0.201'
$ws.Range("E10").Value = '  +10.25%  '
$ws.Range("D11").Value = '''0.652'
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '''54.92'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '''0.0000311'
$ws.Range("E13").Value = '  +4.07%  '
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").Value = '4.145.66'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '''19.50'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '70.645.73'
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '3.573.67'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '''12.52'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '''564.95'
$ws.Range("E20").Value = '  +16.40%  '
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("E23").Value = '  -7.83%  '
$ws.Range("D24").Value = '''4.71'
$ws.Range("E24").Value = '  +7.93%  '
$ws.Range("D25").Value = '''4.98'
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("D26").Value = '''96.04'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").Value = '''3.01'
$ws.Range("E28").Value = '  +2.10%  '
$ws.Range("D29").Value = '''9.21'
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").Value = '''32.44'
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("D32").Value = '''12.57'
$ws.Range("E32").Value = '  +4.83%  '
$ws.Range("D33").Value = '''65.25'
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = '''0.116'
$ws.Range("E34").Value = '  +1.53%  '
$ws.Range("D35").Value = '''3.29'
$ws.Range("E35").Value = '  +4.21%  '
$ws.Range("D36").Value = '''558.54'
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  +5.58%  '
$ws.Range("D38").Value = '''38.33'
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '0.0₃0778'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("D42").Value = '''3.39'
$ws.Range("E42").Value = '  -2.41%  '
$ws.Range("D43").Value = '3.352.97'
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("D44").Value = '''3.09'
$ws.Range("E44").Value = '  -2.43%  '
$ws.Range("D45").Value = '''3.55'
$ws.Range("E45").Value = '  +4.50%  '
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = '''0.0448'
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("D48").Value = '''9.48'
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("E50").Value = '  +19.31%  '
$ws.Range("E51").Value = '  -0.17%  '
